# issue2.xlsx: switch the active tab to Sheet2, move Sheet1's selection off
# of H23 onto A4, and give Sheet2 a formatted value (0.5 shown with the
# built-in "0" number format) plus basic page setup.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1 keeps its data, but the selection moves to A4. Once Sheet2 is
# activated below, Sheet1 naturally stops being the "tabSelected" sheet.
$ws1.Range("A4").Select()

# Sheet2 gets a new numeric value formatted with the built-in "0" number
# format (numFmtId 1), which creates a new cellXfs entry and assigns it to
# A1.
$ws2.Range("A1").Value = 0.5
$ws2.Range("A1").NumberFormat = "0"

# Basic page setup so a <pageSetup> element is emitted for the sheet.
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Sheet2 becomes the active/selected tab (activeTab="1" in workbook.xml,
# tabSelected="1" on its sheetView).
$ws2.Activate()
